# DESIGN/rules/DesignFirst/Main.xlsx -- "Rules" sheet
#
# 1) Collapse / group the data columns (A:K) that carry explicit custom
#    widths, mirroring the outline "collapsed" state recorded on every
#    <col> entry in the sheet.
# 2) D10 should hold the same numeric value as C10 (100), while keeping
#    its own cell style (s="21") untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Group + collapse the sheet's custom-width columns (A:K) -----------
$cols = $ws.Columns("A:K")
$cols.Group()
$cols.ShowDetail = $false

# --- 2. D10 becomes 100 (matching C10), same cell style as before ---------
$ws.Range("D10").Value = 100
